# Daily attendance processing - reverse the order of names/emails in the
# "Recorded By" column (G) wherever multiple entries are comma-separated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row  # xlUp = -4162
if ($lastRow -lt 1) { $lastRow = 1 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if ($val -isnot [string]) { continue }
    if ($val -notmatch ",") { continue }

    $parts = $val -split ",\s*"
    $reversedParts = $parts[($parts.Count - 1)..0]
    $newVal = [string]::Join(", ", $reversedParts)

    if ($newVal -ne $val) {
        $cell.Value2 = $newVal
    }
}
